$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.335.57'
$ws.Range("E2").Value = '  -1.19%  '

$ws.Range("D3").Value = '1.590.80'
$ws.Range("E3").Value = '  -0.39%  '

$ws.Range("E4").Value = '  -0.68%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.85'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.54%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.506'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.97%  '

$ws.Range("E7").Value = '  -0.69%  '

$ws.Range("E9").Value = '  -0.40%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.54'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.37%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0844'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.15%  '

$ws.Range("E12").Value = '  -0.42%  '

$ws.Range("D13").Value = '1.635.05'
$ws.Range("E13").Value = '  +2.20%  '

$ws.Range("E14").Value = '  +0.62%  '

$ws.Range("E15").Value = '  -1.11%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.46'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.47%  '

$ws.Range("D17").Value = '26.338.53'
$ws.Range("E17").Value = '  -1.11%  '

$ws.Range("E18").Value = '  -1.55%  '

$ws.Range("E19").Value = '  +4.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '210.82'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.13%  '

$ws.Range("E21").Value = '  -0.61%  '

$ws.Range("E22").Value = '  -0.29%  '

$ws.Range("E23").Value = '  -4.22%  '

$ws.Range("E24").Value = '  -0.23%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '144.78'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.62%  '

$ws.Range("E26").Value = '  -0.68%  '

$ws.Range("E27").Value = '  -1.43%  '

$ws.Range("E28").Value = '  -0.49%  '

$ws.Range("E29").Value = '  -0.04%  '

$ws.Range("E30").Value = '  -0.21%  '

$ws.Range("E31").Value = '  -0.55%  '

$ws.Range("E32").Value = '  -1.01%  '

$ws.Range("E33").Value = '  +0.70%  '

$ws.Range("D34").Value = '1.308.59'
$ws.Range("E34").Value = '  +2.29%  '

$ws.Range("E35").Value = '  +2.56%  '

$ws.Range("E36").Value = '  -2.13%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.48'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.51%  '

$ws.Range("E38").Value = '  +0.30%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.10'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -13.27%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.811'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.42%  '

$ws.Range("E41").Value = '  -0.56%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.63'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.80%  '

$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.13'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.40%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '62.55'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.02%  '

$ws.Range("B45").Value = 'TrustWalletToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.764'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.62%  '

$ws.Range("D46").Value = '1.726.03'
$ws.Range("E46").Value = '  -0.41%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '88.02'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.53%  '

$ws.Range("E48").Value = '  -1.03%  '

$ws.Range("E49").Value = '  -4.20%  '

$ws.Range("E50").Value = '  -4.35%  '

$ws.Range("E51").Value = '  -1.46%  '
